$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force number-looking Price cells to remain plain text (matches source formatting)
$numericTextCells = @("D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $numericTextCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '19.962.35' ; $ws.Range("E2").Value = '  -8.18%  '
$ws.Range("D3").Value = '1.414.59' ; $ws.Range("E3").Value = '  -8.10%  '
$ws.Range("D4").Value = '1.003' ; $ws.Range("E4").Value = '  +0.18%  '
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '273.33' ; $ws.Range("E6").Value = '  -5.73%  '
$ws.Range("D7").Value = '0.3687' ; $ws.Range("E7").Value = '  -5.56%  '
$ws.Range("D8").Value = '0.3072' ; $ws.Range("E8").Value = '  -3.65%  '
$ws.Range("D9").Value = '39.50' ; $ws.Range("E9").Value = '  -8.06%  '
$ws.Range("D10").Value = '1.001' ; $ws.Range("E10").Value = '  -5.52%  '
$ws.Range("D11").Value = '0.06562' ; $ws.Range("E11").Value = '  -8.94%  '
$ws.Range("D12").Value = '1.001' ; $ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '5.421' ; $ws.Range("E13").Value = '  -4.00%  '
$ws.Range("D14").Value = '17.02' ; $ws.Range("E14").Value = '  -8.74%  '
$ws.Range("D15").Value = '6.165' ; $ws.Range("E15").Value = '  -6.84%  '
$ws.Range("D16").Value = '1.422.88' ; $ws.Range("E16").Value = '  -7.70%  '
$ws.Range("D17").Value = '0.00001004' ; $ws.Range("E17").Value = '  -9.53%  '
$ws.Range("D18").Value = '0.05745' ; $ws.Range("E18").Value = '  -12.79%  '
$ws.Range("D19").Value = '74.05' ; $ws.Range("E19").Value = '  -11.06%  '
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '5.599' ; $ws.Range("E21").Value = '  -9.07%  '
$ws.Range("D22").Value = '14.45' ; $ws.Range("E22").Value = '  -6.25%  '
$ws.Range("D23").Value = '10.91' ; $ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = '2.334'
$ws.Range("D25").Value = '19.979.62' ; $ws.Range("E25").Value = '  -8.14%  '
$ws.Range("D26").Value = '2.274' ; $ws.Range("E26").Value = '  -4.36%  '
$ws.Range("D27").Value = '139.41' ; $ws.Range("E27").Value = '  -4.88%  '
$ws.Range("D28").Value = '16.94' ; $ws.Range("E28").Value = '  -8.01%  '
$ws.Range("D29").Value = '1.580.03' ; $ws.Range("E29").Value = '  -7.86%  '
$ws.Range("D30").Value = '108.94' ; $ws.Range("E30").Value = '  -7.44%  '
$ws.Range("D31").Value = '3.907' ; $ws.Range("E31").Value = '  -19.23%  '
$ws.Range("D32").Value = '5.367' ; $ws.Range("E32").Value = '  -9.41%  '
$ws.Range("D33").Value = '0.8533' ; $ws.Range("E33").Value = '  -12.65%  '
$ws.Range("D34").Value = '0.07714' ; $ws.Range("E34").Value = '  -5.93%  '
$ws.Range("D35").Value = '8.406' ; $ws.Range("E35").Value = '  -4.84%  '
$ws.Range("D36").Value = '0.05736' ; $ws.Range("E36").Value = '  -5.87%  '
$ws.Range("B37").Value = 'Frax' ; $ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' ; $ws.Range("D37").Value = '1.001' ; $ws.Range("E37").Value = '  +0.14%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)' ; $ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' ; $ws.Range("D38").Value = '4.767' ; $ws.Range("E38").Value = '  -7.29%  '
$ws.Range("D39").Value = '10.71' ; $ws.Range("E39").Value = '  +0.13%  '
$ws.Range("B40").Value = 'WEMIXTOKEN' ; $ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' ; $ws.Range("D40").Value = '1.379' ; $ws.Range("E40").Value = '  -6.98%  '
$ws.Range("B41").Value = 'Algorand' ; $ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' ; $ws.Range("D41").Value = '0.1920' ; $ws.Range("E41").Value = '  -5.87%  '
$ws.Range("B42").Value = 'VeChain' ; $ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' ; $ws.Range("D42").Value = '0.02031' ; $ws.Range("E42").Value = '  -7.84%  '
$ws.Range("B43").Value = 'TrustWalletToken' ; $ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' ; $ws.Range("D43").Value = '1.061' ; $ws.Range("E43").Value = '  -10.73%  '
$ws.Range("D44").Value = '0.5293' ; $ws.Range("E44").Value = '  -8.02%  '
$ws.Range("D45").Value = '3.531' ; $ws.Range("E45").Value = '  -5.66%  '
$ws.Range("D46").Value = '12.22' ; $ws.Range("E46").Value = '  -6.81%  '
$ws.Range("D47").Value = '0.5118' ; $ws.Range("E47").Value = '  -7.17%  '
$ws.Range("D48").Value = '1.801' ; $ws.Range("E48").Value = '  -3.83%  '
$ws.Range("D49").Value = '109.26' ; $ws.Range("E49").Value = '  -6.84%  '
$ws.Range("D50").Value = '1.047'
